$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.Execute("Azad Garajev", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $rng.Collapse(0)
    $rng.InsertAfter(" (left 11:18)")
}
